$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of data (row 94) following the same pattern as existing rows
$row = 94
$ws.Cells.Item($row, 1).Value = 46043
$ws.Cells.Item($row, 2).Value = 218
$ws.Cells.Item($row, 3).Value = 224
$ws.Cells.Item($row, 4).Value = 213

# Match the style used by the rest of column A (date-formatted numbers)
$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat
